$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / rename existing header cells (columns already present)
$ws.Range("E3").Value = "Driver"
$ws.Range("F3").Value = "Crew"
$ws.Range("G3").Value = "Truck"
$ws.Range("H3").Value = "Site"
$ws.Range("I3").Value = "Lat"
$ws.Range("J3").Value = "Long"
$ws.Range("L3").Value = "Year Collection"
$ws.Range("M3").Value = "Program"
$ws.Range("N3").Value = "Tank"
$ws.Range("O3").Value = "Trough"
$ws.Range("P3").Value = "Release Method"
$ws.Range("Q3").Value = "Truck Temp"

# New trailing header cells
$ws.Range("R3").Value = "River Temp"
$ws.Range("S3").Value = "Acclimation Time (mins)"
$ws.Range("T3").Value = "Lifestage"
$ws.Range("U3").Value = "Len (cm)"
$ws.Range("V3").Value = "Weight (Kg)"
$ws.Range("W3").Value = "NFish"
$ws.Range("X3").Value = "Comments"

# Give the new cells the same header formatting (grey fill + bottom border) as the rest of the row
$ws.Range("A3").Copy()
$ws.Range("R3:X3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to match the new used range of the header row
[void]$ws.Range("A3:XFD3").Select()
